$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "63.17") need to be
# force-written as text, otherwise Excel auto-converts them to a numeric
# value. We temporarily switch the cell to a text number format, write the
# value, then restore the cell's original style so no formatting changes leak
# into the saved workbook.
function Set-TextValue {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "27.692.94"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.585.45"
$ws.Range("E3").Value = "  -3.09%  "
$ws.Range("E4").Value = "  +0.20%  "
Set-TextValue $ws.Range("D5") "206.75"
$ws.Range("E5").Value = "  -2.35%  "
$ws.Range("E6").Value = "  -3.22%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -4.66%  "
Set-TextValue $ws.Range("D9") "0.254"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "1.810.93"
$ws.Range("E12").Value = "  -3.05%  "
$ws.Range("D13").Value = "1.559.82"
$ws.Range("E13").Value = "  -4.71%  "
$ws.Range("E14").Value = "  -3.73%  "
$ws.Range("E15").Value = "  -5.70%  "
$ws.Range("D16").Value = "27.647.27"
Set-TextValue $ws.Range("D17") "63.17"
$ws.Range("E17").Value = "  -3.37%  "
Set-TextValue $ws.Range("D18") "219.17"
$ws.Range("E18").Value = "  -4.19%  "
$ws.Range("E19").Value = "  -3.54%  "
Set-TextValue $ws.Range("D20") "7.32"
$ws.Range("E20").Value = "  -5.19%  "
$ws.Range("E22").Value = "  -4.95%  "
$ws.Range("E23").Value = "  -5.52%  "
Set-TextValue $ws.Range("D25") "153.75"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("E26").Value = "  +0.25%  "
Set-TextValue $ws.Range("D27") "6.74"
$ws.Range("E27").Value = "  -2.39%  "
Set-TextValue $ws.Range("D28") "15.11"
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("D33").Value = "1.383.70"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("E34").Value = "  -4.79%  "
$ws.Range("E35").Value = "  -5.39%  "
Set-TextValue $ws.Range("D36") "0.964"
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("E39").Value = "  -3.12%  "
Set-TextValue $ws.Range("D40") "0.821"
$ws.Range("E40").Value = "  -3.62%  "
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("E44").Value = "  -3.66%  "
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("E46").Value = "  -4.05%  "
$ws.Range("D47").Value = "1.721.37"
$ws.Range("E47").Value = "  -3.10%  "
Set-TextValue $ws.Range("D48") "88.29"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("E49").Value = "  -1.47%  "
Set-TextValue $ws.Range("D50") "0.0976"
$ws.Range("E50").Value = "  -5.09%  "
$ws.Range("E51").Value = "  -1.09%  "
